# edit.ps1 - Excel COM-interop script
#
# Applies the "cryptos list" price/volume(1h) refresh described by the
# commit "Updated cryptos list on Thu Feb 22 02:24:18 UTC 2024 with GitHub
# Actions": per-coin Price (column D) and Volume(1h) (column E) values are
# refreshed, and two coin pairs (Dai/LEO at rows 27-28, Stellar/ARBITRUM at
# rows 42-43) swap ranking order along with their refreshed data.
#
# Column D holds prices as free-form text (e.g. "51.382.57", "0.999") -
# NOT numbers (some aren't even valid numbers, e.g. thousand-dot-separated
# BTC/ETH prices). A leading apostrophe is used for the numeric-looking
# ones so Excel stores them as text instead of silently parsing them into
# floating point numbers (which would also mangle values like "1.00").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "51.382.57"
$ws.Range("E2").Value = "  -1.40%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.923.87"
$ws.Range("E3").Value = "  -2.54%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.14%  "

# Row 5: BNB
$ws.Range("D5").Value = "'378.09"
$ws.Range("E5").Value = "  +6.75%  "

# Row 6: Solana
$ws.Range("D6").Value = "'103.10"
$ws.Range("E6").Value = "  -3.49%  "

# Row 7: XRP
$ws.Range("D7").Value = "'0.543"
$ws.Range("E7").Value = "  -2.66%  "

# Row 8: USDC
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.08%  "

# Row 9: Cardano
$ws.Range("D9").Value = "'0.586"
$ws.Range("E9").Value = "  -3.86%  "

# Row 10: Avalanche
$ws.Range("D10").Value = "'37.08"
$ws.Range("E10").Value = "  -2.82%  "

# Row 11: TRON
$ws.Range("E11").Value = "  +0.00%  "

# Row 12: Dogecoin
$ws.Range("D12").Value = "'0.0836"
$ws.Range("E12").Value = "  -2.24%  "

# Row 13: Chainlink
$ws.Range("D13").Value = "'18.33"
$ws.Range("E13").Value = "  -3.68%  "

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.382.16"
$ws.Range("E14").Value = "  -2.49%  "

# Row 15: Polkadot
$ws.Range("D15").Value = "'7.35"
$ws.Range("E15").Value = "  -3.63%  "

# Row 16: WrappedEther
$ws.Range("D16").Value = "2.914.82"
$ws.Range("E16").Value = "  -2.69%  "

# Row 17: Polygon
$ws.Range("D17").Value = "'0.930"
$ws.Range("E17").Value = "  -8.51%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "51.309.36"
$ws.Range("E18").Value = "  -1.58%  "

# Row 19: ImmutableX
$ws.Range("D19").Value = "'3.39"
$ws.Range("E19").Value = "  -1.22%  "

# Row 20: Uniswap
$ws.Range("D20").Value = "'7.38"
$ws.Range("E20").Value = "  -1.50%  "

# Row 21: InternetComputer(DFINITY)
$ws.Range("D21").Value = "'12.96"
$ws.Range("E21").Value = "  -4.21%  "

# Row 22: ShibaInu
$ws.Range("D22").Value = "0.0₃0945"
$ws.Range("E22").Value = "  -2.48%  "

# Row 23: Litecoin
$ws.Range("D23").Value = "'68.32"
$ws.Range("E23").Value = "  -1.17%  "

# Row 24: BitcoinCash
$ws.Range("D24").Value = "'261.37"
$ws.Range("E24").Value = "  -1.07%  "

# Row 25: PancakeSwap
$ws.Range("D25").Value = "'2.76"
$ws.Range("E25").Value = "  +1.55%  "

# Row 26: Kaspa
$ws.Range("D26").Value = "'0.170"
$ws.Range("E26").Value = "  -4.19%  "

# Row 27: LEO
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "'4.12"
$ws.Range("E27").Value = "  -5.14%  "

# Row 28: Dai
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.01%  "

# Row 29: EthereumClassic
$ws.Range("D29").Value = "'25.69"
$ws.Range("E29").Value = "  -4.27%  "

# Row 30: Filecoin
$ws.Range("D30").Value = "'7.18"
$ws.Range("E30").Value = "  -3.87%  "

# Row 31: RenderToken
$ws.Range("D31").Value = "'6.85"
$ws.Range("E31").Value = "  +6.96%  "

# Row 32: Hedera
$ws.Range("E32").Value = "  -4.86%  "

# Row 33: Cosmos
$ws.Range("D33").Value = "'9.79"
$ws.Range("E33").Value = "  -4.22%  "

# Row 34: Toncoin
$ws.Range("E34").Value = "  -3.55%  "

# Row 35: OKB
$ws.Range("D35").Value = "'51.08"
$ws.Range("E35").Value = "  -0.04%  "

# Row 36: InjectiveProtocol
$ws.Range("D36").Value = "'34.07"
$ws.Range("E36").Value = "  -5.08%  "

# Row 38: VeChain
$ws.Range("D38").Value = "'0.0421"
$ws.Range("E38").Value = "  -3.82%  "

# Row 39: LidoDAOToken
$ws.Range("D39").Value = "'3.02"
$ws.Range("E39").Value = "  -8.90%  "

# Row 40: Celestia
$ws.Range("D40").Value = "'16.92"
$ws.Range("E40").Value = "  -3.94%  "

# Row 41: Stacks
$ws.Range("D41").Value = "'2.56"

# Row 42: ARBITRUM
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'1.81"
$ws.Range("E42").Value = "  -8.20%  "

# Row 43: Stellar
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "'0.114"
$ws.Range("E43").Value = "  -2.43%  "

# Row 44: Monero
$ws.Range("D44").Value = "'122.90"
$ws.Range("E44").Value = "  -1.43%  "

# Row 45: EnergySwap
$ws.Range("D45").Value = "'21.60"
$ws.Range("E45").Value = "  -5.77%  "

# Row 46: WEMIXToken
$ws.Range("E46").Value = "  -2.83%  "

# Row 47: TheGraph
$ws.Range("D47").Value = "'0.271"
$ws.Range("E47").Value = "  +11.16%  "

# Row 48: Maker
$ws.Range("D48").Value = "2.025.92"
$ws.Range("E48").Value = "  -4.43%  "

# Row 49: ApeXProtocol
$ws.Range("E49").Value = "  -2.39%  "

# Row 50: NEARProtocol
$ws.Range("D50").Value = "'3.16"
$ws.Range("E50").Value = "  -5.10%  "

# Row 51: RocketPoolETH
$ws.Range("D51").Value = "3.202.13"
$ws.Range("E51").Value = "  -2.80%  "
